# Refresh the crypto price/volume snapshot on Sheet1 (coinranking.com pull).
# Columns: A=rank B=Coin C=Link D=Price E=Volume(1h). Rows 2-51 are data rows;
# row 16 (Dai) is unchanged this run, and a handful of rows only move on
# Volume(1h) because Price held steady between refreshes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values keyed by row number. $null Price means "leave column D alone".
$updates = @(
    @{ Row = 2; Price = "25.570.93"; Volume = "  +2.28%  " },
    @{ Row = 3; Price = "1.669.27"; Volume = "  +1.84%  " },
    @{ Row = 4; Price = "0.9983"; Volume = "  +0.15%  " },
    @{ Row = 5; Price = "235.09"; Volume = "  +0.87%  " },
    @{ Row = 6; Price = "0.9994"; Volume = "  +0.04%  " },
    @{ Row = 7; Price = "0.4648"; Volume = "  -3.17%  " },
    @{ Row = 8; Price = "0.2576"; Volume = "  -1.02%  " },
    @{ Row = 9; Price = "0.06132"; Volume = "  +0.51%  " },
    @{ Row = 10; Price = "1.664.59"; Volume = "  +1.61%  " },
    @{ Row = 11; Price = "0.06950"; Volume = "  -1.55%  " },
    @{ Row = 12; Price = "14.65"; Volume = "  +0.77%  " },
    @{ Row = 13; Price = "4.341"; Volume = "  -0.85%  " },
    @{ Row = 14; Price = "74.91"; Volume = "  +1.92%  " },
    @{ Row = 15; Price = "0.5720"; Volume = "  -4.53%  " },
    @{ Row = 17; Price = "0.9995"; Volume = "  +0.14%  " },
    @{ Row = 18; Price = "25.562.46"; Volume = "  +2.29%  " },
    @{ Row = 19; Price = "0.000006735"; Volume = "  +2.08%  " },
    @{ Row = 20; Price = $null; Volume = "  +0.89%  " },
    @{ Row = 21; Price = "1.878.32"; Volume = "  +1.66%  " },
    @{ Row = 22; Price = "4.413"; Volume = "  +0.69%  " },
    @{ Row = 23; Price = "8.670"; Volume = "  +0.90%  " },
    @{ Row = 24; Price = "5.233"; Volume = "  -0.40%  " },
    @{ Row = 25; Price = "134.38"; Volume = "  +0.87%  " },
    @{ Row = 26; Price = "14.86"; Volume = "  -0.34%  " },
    @{ Row = 27; Price = "1.364"; Volume = "  -1.43%  " },
    @{ Row = 28; Price = "1.713"; Volume = "  +4.19%  " },
    @{ Row = 29; Price = "103.82"; Volume = "  -0.45%  " },
    @{ Row = 30; Price = "3.959"; Volume = "  +2.37%  " },
    @{ Row = 31; Price = "0.07701"; Volume = "  -0.06%  " },
    @{ Row = 32; Price = "3.598"; Volume = "  +1.33%  " },
    @{ Row = 33; Price = "0.04316"; Volume = "  +0.80%  " },
    @{ Row = 34; Price = "2.619"; Volume = "  +1.81%  " },
    @{ Row = 35; Price = "0.9440"; Volume = "  +1.82%  " },
    @{ Row = 36; Price = "0.6000"; Volume = "  +2.34%  " },
    @{ Row = 37; Price = "0.9182"; Volume = "  +10.97%  " },
    @{ Row = 38; Price = "2.480"; Volume = "  -2.83%  " },
    @{ Row = 39; Price = "105.10"; Volume = "  +6.79%  " },
    @{ Row = 40; Price = "0.9991"; Volume = "  +0.10%  " },
    @{ Row = 41; Price = "1.843"; Volume = "  +5.39%  " },
    @{ Row = 42; Price = "0.01463"; Volume = "  -3.86%  " },
    @{ Row = 43; Price = "5.041"; Volume = "  +7.60%  " },
    @{ Row = 44; Price = $null; Volume = "  +0.26%  " },
    @{ Row = 45; Price = "0.1110"; Volume = "  +2.08%  " },
    @{ Row = 46; Price = "0.05251"; Volume = "  +1.10%  " },
    @{ Row = 47; Price = "6.124"; Volume = "  +0.77%  " },
    @{ Row = 48; Price = "29.73"; Volume = "  +1.80%  " },
    @{ Row = 49; Price = "7.560"; Volume = "  +4.67%  " },
    @{ Row = 50; Price = $null; Volume = "  +0.32%  " },
    @{ Row = 51; Price = $null; Volume = "  +0.25%  " }

)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D$($u.Row)")
        # Price strings like "0.9983" or "14.65" parse as plain numbers, which
        # would make Excel silently coerce them to numeric cells (dropping the
        # fixed-width formatting the feed relies on, e.g. "0.9990" -> 0.999).
        # Values with two dots (e.g. "25.570.93") can never parse as a number,
        # so they round-trip as text without any extra help. Force the rest to
        # text with a leading apostrophe, then strip the resulting quote-prefix
        # style back to Normal so formatting matches the untouched cells.
        if ($u.Price -match '^-?[0-9]*\.?[0-9]+$') {
            $priceCell.Value = "'" + $u.Price
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $u.Price
        }
    }
    $ws.Range("E$($u.Row)").Value = $u.Volume
}
